$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70-77 down to 71-78
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the weekly data point
$ws.Cells.Item(70, 1).Value = 8
$ws.Cells.Item(70, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44449
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 100112001
$ws.Cells.Item(70, 7).Value = "Berenjena"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 600
$ws.Cells.Item(70, 11).Value = 9000
$ws.Cells.Item(70, 12).Value = 10000
$ws.Cells.Item(70, 13).Value = 9500
$ws.Cells.Item(70, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(70, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value = 158
$ws.Cells.Item(70, 17).Value = 60
$ws.Cells.Item(70, 18).Value = "Hortaliza"
